# Scheduled market-data refresh: push newly-fetched Universalis price
# snapshots (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ)
# and the resulting profit deltas into each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 8006038
$ws.Range("I5").Value = 5007548
$ws.Range("K5").Value = 5007548
$ws.Range("M5").Value = -5007433

$ws.Range("H76").Value = 7777
$ws.Range("I76").Value = 7777
$ws.Range("K76").Value = 7777
$ws.Range("M76").Value = -7462

$ws.Range("H79").Value = 7777
$ws.Range("I79").Value = 7777
$ws.Range("K79").Value = 7777
$ws.Range("M79").Value = -6685

$ws.Range("H86").Value = 4333.3335
$ws.Range("J86").Value = 3000
$ws.Range("L86").Value = 3000
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 4333.3335
$ws.Range("J89").Value = 3000
$ws.Range("L89").Value = 15000
$ws.Range("N89").Value = -26232

$ws.Range("H92").Value = 1055
$ws.Range("I92").Value = 450
$ws.Range("K92").Value = 450
$ws.Range("M92").Value = 798

$ws.Range("H116").Value = 5048.5
$ws.Range("I116").Value = 5368.5
$ws.Range("J116").Value = 4248.5
$ws.Range("K116").Value = 5368.5
$ws.Range("L116").Value = 4248.5
$ws.Range("M116").Value = -1926.5
$ws.Range("N116").Value = -11132.5

$ws.Range("H138").Value = 5284.55
$ws.Range("J138").Value = 5678.1143
$ws.Range("L138").Value = 17034.3429
$ws.Range("N138").Value = -27314.3429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H32").Value = 2676.3635
$ws.Range("I32").Value = 1918.0704
$ws.Range("K32").Value = 1918.0704
$ws.Range("M32").Value = -1631.0704

$ws.Range("H63").Value = 7849.4
$ws.Range("I63").Value = 7415.8335
$ws.Range("J63").Value = 8499.75
$ws.Range("K63").Value = 7415.8335
$ws.Range("L63").Value = 8499.75
$ws.Range("M63").Value = -6729.8335
$ws.Range("N63").Value = -9871.75

$ws.Range("H66").Value = 7849.4
$ws.Range("I66").Value = 7415.8335
$ws.Range("J66").Value = 8499.75
$ws.Range("K66").Value = 37079.1675
$ws.Range("L66").Value = 42498.75
$ws.Range("M66").Value = -33647.1675
$ws.Range("N66").Value = -49362.75

$ws.Range("H118").Value = 110000
$ws.Range("J118").Value = 110000
$ws.Range("L118").Value = 110000
$ws.Range("N118").Value = -113314

$ws.Range("H132").Value = 3653.5264
$ws.Range("I132").Value = 2439.5454
$ws.Range("J132").Value = 5322.75
$ws.Range("K132").Value = 7318.6362
$ws.Range("L132").Value = 15968.25
$ws.Range("M132").Value = -4788.6362
$ws.Range("N132").Value = -21028.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 4716.5
$ws.Range("I36").Value = 4716.5
$ws.Range("K36").Value = 4716.5
$ws.Range("M36").Value = -4182.5

$ws.Range("H39").Value = 1418.75
$ws.Range("I39").Value = 1380
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 1380
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = -991
$ws.Range("N39").Value = -2778

$ws.Range("H75").Value = 19999
$ws.Range("I75").Value = 19999
$ws.Range("K75").Value = 19999
$ws.Range("M75").Value = -19063

$ws.Range("H78").Value = 19999
$ws.Range("I78").Value = 19999
$ws.Range("K78").Value = 59997
$ws.Range("M78").Value = -55317

$ws.Range("H99").Value = 1451.9231
$ws.Range("I99").Value = 1443.7273
$ws.Range("J99").Value = 1497
$ws.Range("K99").Value = 1443.7273
$ws.Range("L99").Value = 1497
$ws.Range("M99").Value = 54.27269999999999
$ws.Range("N99").Value = -4493

$ws.Range("H134").Value = 2706
$ws.Range("I134").Value = 2532.842
$ws.Range("J134").Value = 3176
$ws.Range("K134").Value = 7598.526
$ws.Range("L134").Value = 9528
$ws.Range("M134").Value = -5063.526
$ws.Range("N134").Value = -14598

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8182545
$ws.Range("J6").Value = 999
$ws.Range("L6").Value = 999
$ws.Range("N6").Value = -1225

$ws.Range("H62").Value = 8681
$ws.Range("I62").Value = 8211.4
$ws.Range("J62").Value = 12203
$ws.Range("K62").Value = 8211.4
$ws.Range("L62").Value = 12203
$ws.Range("M62").Value = -7587.4
$ws.Range("N62").Value = -13451

$ws.Range("H65").Value = 8681
$ws.Range("I65").Value = 8211.4
$ws.Range("J65").Value = 12203
$ws.Range("K65").Value = 41057
$ws.Range("L65").Value = 61015
$ws.Range("M65").Value = -37937
$ws.Range("N65").Value = -67255

$ws.Range("H99").Value = 25576
$ws.Range("I99").Value = 23488.223
$ws.Range("J99").Value = 27924.75
$ws.Range("K99").Value = 23488.223
$ws.Range("L99").Value = 27924.75
$ws.Range("M99").Value = -21990.223
$ws.Range("N99").Value = -30920.75

$ws.Range("H107").Value = 406.07407
$ws.Range("I107").Value = 393.83334
$ws.Range("J107").Value = 504
$ws.Range("K107").Value = 393.83334
$ws.Range("L107").Value = 504
$ws.Range("M107").Value = 1526.16666
$ws.Range("N107").Value = -4344

$ws.Range("H126").Value = 25576
$ws.Range("I126").Value = 23488.223
$ws.Range("J126").Value = 27924.75
$ws.Range("K126").Value = 70464.66900000001
$ws.Range("L126").Value = 83774.25
$ws.Range("M126").Value = -67994.66900000001
$ws.Range("N126").Value = -88714.25

$ws.Range("H132").Value = 3938.9333
$ws.Range("I132").Value = 3242.889
$ws.Range("J132").Value = 4983
$ws.Range("K132").Value = 9728.667000000001
$ws.Range("L132").Value = 14949
$ws.Range("M132").Value = -7198.667000000001
$ws.Range("N132").Value = -20009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1228
$ws.Range("I12").Value = 1001.55554
$ws.Range("J12").Value = 1397.8334
$ws.Range("K12").Value = 3004.66662
$ws.Range("L12").Value = 4193.5002
$ws.Range("M12").Value = -2831.66662
$ws.Range("N12").Value = -4539.5002

$ws.Range("H121").Value = 752.5
$ws.Range("I121").Value = 893.3333
$ws.Range("K121").Value = 2679.9999
$ws.Range("M121").Value = -1369.9999

$ws.Range("H132").Value = 1665.75
$ws.Range("I132").Value = 1665
$ws.Range("J132").Value = 1666
$ws.Range("K132").Value = 14985
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -12455
$ws.Range("N132").Value = -20054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2019.9
$ws.Range("I16").Value = 2144.3333
$ws.Range("K16").Value = 2144.3333
$ws.Range("M16").Value = -1974.3333

$ws.Range("H32").Value = 1859.8
$ws.Range("I32").Value = 1859.8
$ws.Range("K32").Value = 1859.8
$ws.Range("M32").Value = -1542.8

$ws.Range("H93").Value = 2717.8
$ws.Range("I93").Value = 2147.25
$ws.Range("K93").Value = 2147.25
$ws.Range("M93").Value = -899.25

$ws.Range("H136").Value = 5259.04
$ws.Range("I136").Value = 4325.1763
$ws.Range("K136").Value = 12975.5289
$ws.Range("M136").Value = -10425.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 20936.2
$ws.Range("J45").Value = 20936.2
$ws.Range("L45").Value = 20936.2
$ws.Range("N45").Value = -21918.2

$ws.Range("H74").Value = 19999.5
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 19999.5
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
